# Add a new "09_24" build entry to the bottom of the test-log sheet and
# document the level-7/8 difficulty increase + bug fixes (per commit msg).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 33 is the current last data row; the new entry goes to row 34.
$srcRow = 33
$newRow = 34

# Copy formatting (styles 3/3/3/5 -> vertical-top [+wrap on E]) from the
# previous row so the new row renders identically to its neighbours.
$ws.Range("B$srcRow`:E$srcRow").Copy()
$ws.Range("B$newRow`:E$newRow").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Version / Testmethoden / Tester / Neue Features (columns B..E)
$ws.Range("B$newRow").Value = "DiscordiaAgency_Demo_2017_09_24.exe"
$ws.Range("C$newRow").Value = "Entwicklung"
$ws.Range("D$newRow").Value = "Anna Franziska"
$ws.Range("E$newRow").Value = "Hauptmenü: keine Fehlermeldung mehr, dass Objekt nicht gefunden; Hauptmenü: Musik startet nicht mehr erneut, wenn zur Steuerungsübersicht gewechselt wird; zufälliges Rotieren der Wachen funktioniert jetzt auch um 0 Grad herum; Level 7 und 8 schwerer gemacht; Verkleiden klappt, auch wenn man über 2 toten Körpern steht"

# Match the auto-fit row height Excel would compute for the wrapped text.
$ws.Rows($newRow).RowHeight = 105

# Point the selection at the newly added row, same as the source workbook
# (the sheet's frozen header pane is left untouched).
[void]$ws.Range("E$newRow").Select()
